$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM data changes ("Change resistor values back to 1k") ---

# Row 29: R8;R9;R10 group grows to include R11-R18 (all now 1k), QTY 3 -> 11
$ws.Cells.Item(29, 1).Value = "R8;R9;R10;R11;R12;R13;R14;R15;R16;R17;R18"
$ws.Cells.Item(29, 4).Value = 11

# Rows 35/36: R3 (3k8) now listed before R4 (3r9); R4's value becomes lowercase "3r9"
$ws.Cells.Item(35, 1).Value = "R3"
$ws.Cells.Item(35, 2).Value = "3k8"
$ws.Cells.Item(35, 3).Value = "Resistors_SMD:R_0402"
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(35, 5).Value = "Susumu"
$ws.Cells.Item(35, 6).Value = "RR0510P-3831-D"
$ws.Cells.Item(35, 7).Value = "Mouser"
$ws.Cells.Item(35, 8).Value = "754-RR0510P-3831D"

$ws.Cells.Item(36, 1).Value = "R4"
$ws.Cells.Item(36, 2).Value = "3r9"
$ws.Cells.Item(36, 3).Value = "Resistors_SMD:R_0402"
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(36, 5).Value = "Vishay"
$ws.Cells.Item(36, 6).Value = "CRCW04023R90FKED"
$ws.Cells.Item(36, 7).Value = "Mouser"
$ws.Cells.Item(36, 8).Value = "71-CRCW04023R90FKED"

# Row 38 (R11;R12;...;R18 @ 75) is gone now (folded into row 29); PI1 row shifts up from 39 to 38
$ws.Rows(38).Delete()

# --- Cosmetic / view changes ---

# Zoom 100% -> 110%
$excel.ActiveWindow.Zoom = 110

# Column A got wider to fit the longer Refs text
$ws.Columns.Item(1).ColumnWidth = 39.5

# Print paper size "A4"(9) -> "Letter"(1)
$ws.PageSetup.PaperSize = 1

# Header/footer font label + footer wording changed
$ws.PageSetup.CenterHeader = "&`"Times New Roman,Regular`"&12&A"
$ws.PageSetup.CenterFooter = "&`"Times New Roman,Regular`"&12Page &P"
